# Update countries & provincias Spain
# Applies the COVID data refresh: updates case numbers for several
# countries and reflects the resulting ranking swaps (Liberia overtakes
# Santo Tome y Principe; Togo overtakes Birmania and Suazilandia).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Straightforward numeric updates (no row-order change) ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 1375129
$ws.Range("C4").Value = 7491
$ws.Range("D4").Value = 258037
$ws.Range("E4").Value = 1035993
$ws.Range("F4").Value = 16473
$ws.Range("G4").Value = 312
$ws.Range("H4").Value = 81099

# Row 15: India
$ws.Range("B15").Value = 70765
$ws.Range("C15").Value = 3604
$ws.Range("E15").Value = 45922
$ws.Range("G15").Value = 82
$ws.Range("H15").Value = 2294

# Row 31: Irlanda
$ws.Range("B31").Value = 23135
$ws.Range("C31").Value = 139
$ws.Range("E31").Value = 4558
$ws.Range("G31").Value = 9
$ws.Range("H31").Value = 1467

# Row 100: El Salvador
$ws.Range("E100").Value = 615
$ws.Range("G100").Value = 1
$ws.Range("H100").Value = 18

# Row 121: Jordania
$ws.Range("B121").Value = 562
$ws.Range("C121").Value = 22
$ws.Range("D121").Value = 390
$ws.Range("E121").Value = 163

# Row 181: Eritrea
$ws.Range("D181").Value = 38
$ws.Range("E181").Value = 1

# --- Ranking swaps ---
# Liberia's total cases (199 -> 211) overtakes Santo Tome y Principe's 208,
# so Liberia now occupies row 142 (with its new totals) and Santo Tome y
# Principe slides down to row 143 (keeping its prior totals unchanged).

$ws.Range("A142").Value = "Liberia"
$ws.Range("B142").Value = 211
$ws.Range("C142").Value = 12
$ws.Range("D142").Value = 85
$ws.Range("E142").Value = 106
$ws.Range("F142").Value = 0
$ws.Range("G142").Value = 0
$ws.Range("H142").Value = 20

$ws.Range("A143").Value = "Santo Tome y Principe"
$ws.Range("B143").Value = 208
$ws.Range("C143").Value = 0
$ws.Range("D143").Value = 4
$ws.Range("E143").Value = 199
$ws.Range("F143").Value = 0
$ws.Range("G143").Value = 0
$ws.Range("H143").Value = 5

# Togo's total cases (174 -> 181) overtakes Birmania's 180, so Togo now
# occupies row 148 (with its new totals); Birmania and Suazilandia each
# slide down one row, keeping their prior totals unchanged.

$ws.Range("A148").Value = "Togo"
$ws.Range("B148").Value = 181
$ws.Range("C148").Value = 7
$ws.Range("D148").Value = 89
$ws.Range("E148").Value = 81
$ws.Range("F148").Value = 0
$ws.Range("G148").Value = 0
$ws.Range("H148").Value = 11

$ws.Range("A149").Value = "Birmania"
$ws.Range("B149").Value = 180
$ws.Range("C149").Value = 0
$ws.Range("D149").Value = 74
$ws.Range("E149").Value = 100
$ws.Range("F149").Value = 0
$ws.Range("G149").Value = 0
$ws.Range("H149").Value = 6

$ws.Range("A150").Value = "Suazilandia"
$ws.Range("B150").Value = 175
$ws.Range("C150").Value = 3
$ws.Range("D150").Value = 28
$ws.Range("E150").Value = 145
$ws.Range("F150").Value = 0
$ws.Range("G150").Value = 0
$ws.Range("H150").Value = 2
